# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - match formatting of the existing header row
# (bold/bordered/centered style already applied to A1:AC1) by copying an
# existing header cell's formatting, then overwriting the text.
$headerSrc = $ws.Range("AC1")

$wins = $ws.Range("AD1")
$headerSrc.Copy($wins)
$wins.Value = "Wins"

$losses = $ws.Range("AE1")
$headerSrc.Copy($losses)
$losses.Value = "Losses"

$ties = $ws.Range("AF1")
$headerSrc.Copy($ties)
$ties.Value = "Ties"

# Fill in the team record for every data row (2-49): every row gets the
# same season totals (67 wins, 95 losses, 0 ties).
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 67   # AD
    $ws.Cells.Item($r, 31).Value = 95   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
